# Scoreboard.xlsx update — "Add files via upload"
# - Adds semifinal qualifier results (rows 2-9) to the SFM and SFF sheets
# - Renames one SFF team to "Muscles & Brains"
# - Updates the active sheet/selection state to match the saved workbook view

$wb = $excel.ActiveWorkbook

$wsScoreM = $wb.Worksheets.Item("ScoreM")
$wsScoreF = $wb.Worksheets.Item("ScoreF")
$wsSFM    = $wb.Worksheets.Item("SFM")
$wsSFF    = $wb.Worksheets.Item("SFF")

# ---- SFM (Semifinal Men) results ----
$wsSFM.Range("A2").Value = "Kjetil og Kasper"
$wsSFM.Range("B2").Value = "Hallo"
$wsSFM.Range("C2").Value = 20
$wsSFM.Range("D2").Value = 5
$wsSFM.Range("E2").Value = 37
$wsSFM.Range("F2").Value = 860
$wsSFM.Range("G2").Value = 5
$wsSFM.Range("H2").Value = 0
$wsSFM.Range("I2").Value = 88
$wsSFM.Range("A3").Value = "Laszlo og Marcus"
$wsSFM.Range("B3").Value = "Exorcisers"
$wsSFM.Range("C3").Value = 18
$wsSFM.Range("D3").Value = 6
$wsSFM.Range("E3").Value = 0
$wsSFM.Range("F3").Value = 823
$wsSFM.Range("G3").Value = 5
$wsSFM.Range("H3").Value = 0
$wsSFM.Range("I3").Value = 109
$wsSFM.Range("A4").Value = "Ole og Mikus"
$wsSFM.Range("B4").Value = "Fit and Lazy"
$wsSFM.Range("C4").Value = 16
$wsSFM.Range("D4").Value = 3
$wsSFM.Range("E4").Value = 56
$wsSFM.Range("F4").Value = 860
$wsSFM.Range("G4").Value = 5
$wsSFM.Range("H4").Value = 0
$wsSFM.Range("I4").Value = 108
$wsSFM.Range("A5").Value = "Morten og Gabor"
$wsSFM.Range("B5").Value = "Daddszkys"
$wsSFM.Range("C5").Value = 14
$wsSFM.Range("D5").Value = 5
$wsSFM.Range("E5").Value = 34
$wsSFM.Range("F5").Value = 860
$wsSFM.Range("G5").Value = 5
$wsSFM.Range("H5").Value = 0
$wsSFM.Range("I5").Value = 77
$wsSFM.Range("A6").Value = "Jakob og Finn"
$wsSFM.Range("B6").Value = "The NHH Nerds"
$wsSFM.Range("C6").Value = 12
$wsSFM.Range("D6").Value = 4
$wsSFM.Range("E6").Value = 47
$wsSFM.Range("F6").Value = 860
$wsSFM.Range("G6").Value = 5
$wsSFM.Range("H6").Value = 0
$wsSFM.Range("I6").Value = 110
$wsSFM.Range("A7").Value = "Jonas og Arild"
$wsSFM.Range("B7").Value = "Team Kongobajer"
$wsSFM.Range("C7").Value = 10
$wsSFM.Range("D7").Value = 4
$wsSFM.Range("E7").Value = 26
$wsSFM.Range("F7").Value = 860
$wsSFM.Range("G7").Value = 5
$wsSFM.Range("H7").Value = 0
$wsSFM.Range("I7").Value = 84
$wsSFM.Range("A8").Value = "Eirik og Knut"
$wsSFM.Range("B8").Value = "To privilegerte menn"
$wsSFM.Range("C8").Value = 8
$wsSFM.Range("D8").Value = 3
$wsSFM.Range("E8").Value = 0
$wsSFM.Range("F8").Value = 860
$wsSFM.Range("G8").Value = 5
$wsSFM.Range("H8").Value = 0
$wsSFM.Range("I8").Value = 108
$wsSFM.Range("A9").Value = "Rashad og Emil"
$wsSFM.Range("B9").Value = "RX on the Beach"
$wsSFM.Range("C9").Value = 6
$wsSFM.Range("D9").Value = 5
$wsSFM.Range("E9").Value = 3
$wsSFM.Range("F9").Value = 860
$wsSFM.Range("G9").Value = 5
$wsSFM.Range("H9").Value = 0
$wsSFM.Range("I9").Value = 96

# ---- SFF (Semifinal Women) results ----
$wsSFF.Range("A2").Value = "Anne og Sara"
$wsSFF.Range("B2").Value = "How I met this runner"
$wsSFF.Range("C2").Value = 20
$wsSFF.Range("D2").Value = 5
$wsSFF.Range("E2").Value = 37
$wsSFF.Range("F2").Value = 860
$wsSFF.Range("G2").Value = 5
$wsSFF.Range("H2").Value = 0
$wsSFF.Range("I2").Value = 111
$wsSFF.Range("A3").Value = "Martine og Oda"
$wsSFF.Range("B3").Value = "Team AnabOle"
$wsSFF.Range("C3").Value = 18
$wsSFF.Range("D3").Value = 4
$wsSFF.Range("E3").Value = 28
$wsSFF.Range("F3").Value = 860
$wsSFF.Range("G3").Value = 5
$wsSFF.Range("H3").Value = 0
$wsSFF.Range("I3").Value = 121
$wsSFF.Range("A4").Value = "Linn Therese og Julie"
$wsSFF.Range("B4").Value = "VB Hestejenter"
$wsSFF.Range("C4").Value = 16
$wsSFF.Range("D4").Value = 6
$wsSFF.Range("E4").Value = 0
$wsSFF.Range("F4").Value = 858
$wsSFF.Range("G4").Value = 5
$wsSFF.Range("H4").Value = 0
$wsSFF.Range("I4").Value = 124
$wsSFF.Range("A5").Value = "Helene og Karoline"
$wsSFF.Range("B5").Value = "Show me your jerk"
$wsSFF.Range("C5").Value = 14
$wsSFF.Range("D5").Value = 5
$wsSFF.Range("E5").Value = 36
$wsSFF.Range("F5").Value = 860
$wsSFF.Range("G5").Value = 5
$wsSFF.Range("H5").Value = 0
$wsSFF.Range("I5").Value = 101
$wsSFF.Range("A6").Value = "Monica og Mai"
$wsSFF.Range("B6").Value = "The Ones"
$wsSFF.Range("C6").Value = 12
$wsSFF.Range("D6").Value = 4
$wsSFF.Range("E6").Value = 53
$wsSFF.Range("F6").Value = 860
$wsSFF.Range("G6").Value = 5
$wsSFF.Range("H6").Value = 0
$wsSFF.Range("I6").Value = 86
$wsSFF.Range("A7").Value = "Gøril og Cecilie"
$wsSFF.Range("B7").Value = "Muscles & Brains"
$wsSFF.Range("C7").Value = 10
$wsSFF.Range("D7").Value = 6
$wsSFF.Range("E7").Value = 0
$wsSFF.Range("F7").Value = 811
$wsSFF.Range("G7").Value = 5
$wsSFF.Range("H7").Value = 0
$wsSFF.Range("I7").Value = 119
$wsSFF.Range("A8").Value = "Elise og Kamilla"
$wsSFF.Range("B8").Value = "Regnbuebarna"
$wsSFF.Range("C8").Value = 8
$wsSFF.Range("D8").Value = 6
$wsSFF.Range("E8").Value = 0
$wsSFF.Range("F8").Value = 834
$wsSFF.Range("G8").Value = 5
$wsSFF.Range("H8").Value = 0
$wsSFF.Range("I8").Value = 112
$wsSFF.Range("A9").Value = "Lina og Maren"
$wsSFF.Range("B9").Value = "Comeback Kids"
$wsSFF.Range("C9").Value = 6
$wsSFF.Range("D9").Value = 6
$wsSFF.Range("E9").Value = 0
$wsSFF.Range("F9").Value = 844
$wsSFF.Range("G9").Value = 5
$wsSFF.Range("H9").Value = 0
$wsSFF.Range("I9").Value = 85

# ---- Selections on the sheets that keep their previous layout ----
$wsScoreM.Range("C24").Select()
$wsScoreF.Range("G13").Select()

# ---- SFF keeps its data selected but is no longer the active tab ----
$wsSFF.Range("A2:I9").Select()

# ---- SFM becomes the active sheet/tab (activate last so it "wins") ----
$wsSFM.Activate()
$wsSFM.Range("B7").Select()
